# macro to save each sheet in csv
# (this pass only captures the resulting workbook state: numbering the
#  quote rows on sheet "03" and leaving sheet "09" as the active tab)

$wb = $excel.ActiveWorkbook

# --- Sheet "03": add a running index in column A (rows 2-22), 1..21 ---
$ws03 = $wb.Worksheets.Item("03")

for ($i = 1; $i -le 21; $i++) {
    $row = $i + 1
    $ws03.Cells.Item($row, 1).Value = $i
}

# Update the sheet's selection to the newly filled column A range
$ws03.Range("A2:A22").Select()

# --- Make sheet "09" the active/visible tab, as in the target workbook ---
$ws09 = $wb.Worksheets.Item("09")
$ws09.Activate()

Write-Output "Updated sheet '03' index column and activated sheet '09'."
